# Update "想去人数" (number of people interested) counts on the
# 展览 (Exhibition) and 全部类型 (All types) sheets.
#   展览·龙泉ACG动漫游戏博览会  (row 3): F3  1426 -> 1428
#   展览·第四届HP国风动漫游戏嘉年华 (row 5): F5  10   -> 11

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 1428
    $ws.Range("F5").Value = 11
}
